$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 16 (TS_010 Header Desktop) test-case count: 16 -> 24
$ws.Range("F16").Value = 24

# --- Row 18 first: (TS_010.2) Sub Desktop Mac ---
# (values entered in this order so the shared-string table matches the
# author's original authoring order)
$ws.Range("D18").Value = "Validate the functionality of home page > Header > Desktop> Sub Desktop Mac."
$ws.Range("B18").Value = "(TS_010.2)`nSub Desktop Mac"

# --- Row 17: (TS_010.1) Sub Desktop PC ---
$ws.Range("D17").Value = "Validate the functionality of home page > Header > Desktop> Sub Desktop PC."
$ws.Range("B17").Value = "(TS_010.1)`nSub Desktop PC"

# --- Row 19: (TS_011) Cart ---
$ws.Range("B19").Value = "(TS_011)`n Cart"
$ws.Range("D19").Value = "Validate the functionality of the product cart  across the app."

# Common column values / counts
$ws.Range("C17").Value = "FRS"
$ws.Range("C18").Value = "FRS"
$ws.Range("C19").Value = "FRS"
$ws.Range("F17").Value = 5
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 16

# Copy direct formatting (B/C/D/F styles) from row 16 onto the new rows,
# one column at a time so no stray E-column cell/style gets introduced.
foreach ($r in 17..19) {
    $ws.Range("B16").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)
    $ws.Range("C16").Copy()
    $ws.Range("C$r").PasteSpecial(-4122)
    $ws.Range("D16").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)
    $ws.Range("F16").Copy()
    $ws.Range("F$r").PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = 30
}
$ws.Application.CutCopyMode = $false

# re-apply the values/counts since PasteSpecial(formats) shouldn't disturb
# them, but make sure they're correct regardless of paste ordering
$ws.Range("B17").Value = "(TS_010.1)`nSub Desktop PC"
$ws.Range("C17").Value = "FRS"
$ws.Range("D17").Value = "Validate the functionality of home page > Header > Desktop> Sub Desktop PC."
$ws.Range("F17").Value = 5

$ws.Range("B18").Value = "(TS_010.2)`nSub Desktop Mac"
$ws.Range("C18").Value = "FRS"
$ws.Range("D18").Value = "Validate the functionality of home page > Header > Desktop> Sub Desktop Mac."
$ws.Range("F18").Value = 2

$ws.Range("B19").Value = "(TS_011)`n Cart"
$ws.Range("C19").Value = "FRS"
$ws.Range("D19").Value = "Validate the functionality of the product cart  across the app."
$ws.Range("F19").Value = 16

# Row 23: grand total of test cases
$ws.Range("F23").Formula = "=SUM(F4:F22)"

# Update the sheet view (scrolled position / active selection) to match
# the post-edit state
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("D20").Select()
